$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '66.346.68'
$ws.Range('E2').Value = '  -0.84%  '
$ws.Range('D3').Value = '3.216.05'
$ws.Range('E3').Value = '  +0.33%  '
$ws.Range('E4').Value = '  -0.19%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '605.56'
$ws.Range('E5').Value = '  +0.19%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '155.15'
$ws.Range('E6').Value = '  -1.41%  '
$ws.Range('E7').Value = '  -0.05%  '
$ws.Range('D8').Value = '3.214.76'
$ws.Range('E8').Value = '  +0.32%  '
$ws.Range('E9').Value = '  -2.30%  '
$ws.Range('E10').Value = '  -1.08%  '
$ws.Range('E11').Value = '  -3.68%  '
$ws.Range('E12').Value = '  -3.36%  '
$ws.Range('E13').Value = '  -1.07%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '38.24'
$ws.Range('E14').Value = '  -2.73%  '
$ws.Range('D15').Value = '3.742.53'
$ws.Range('E15').Value = '  +0.08%  '
$ws.Range('D16').Value = '66.396.14'
$ws.Range('E16').Value = '  -0.74%  '
$ws.Range('D17').Value = '3.216.51'
$ws.Range('E17').Value = '  +0.00%  '
$ws.Range('E18').Value = '  -3.43%  '
$ws.Range('E19').Value = '  +0.78%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '506.53'
$ws.Range('E20').Value = '  -3.32%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '15.20'
$ws.Range('E21').Value = '  -1.93%  '
$ws.Range('E22').Value = '  -2.41%  '
$ws.Range('E23').Value = '  -3.25%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '14.50'
$ws.Range('E24').Value = '  -3.99%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '85.01'
$ws.Range('E25').Value = '  -0.68%  '
$ws.Range('B26').Value = 'Hedera'
$ws.Range('C26').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.151'
$ws.Range('E26').Value = '  +66.25%  '
$ws.Range('B27').Value = 'Dai'
$ws.Range('C27').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '1.00'
$ws.Range('E27').Value = '  +0.01%  '
$ws.Range('E28').Value = '  -1.23%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '8.99'
$ws.Range('E29').Value = '  -3.29%  '
$ws.Range('E30').Value = '  -2.13%  '
$ws.Range('B31').Value = 'Stacks'
$ws.Range('C31').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '2.89'
$ws.Range('E31').Value = '  -3.92%  '
$ws.Range('B32').Value = 'NEARProtocol'
$ws.Range('C32').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '6.90'
$ws.Range('E32').Value = '  -1.59%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '28.25'
$ws.Range('E33').Value = '  -0.44%  '
$ws.Range('E34').Value = '  +0.04%  '
$ws.Range('E35').Value = '  -4.91%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '6.38'
$ws.Range('E36').Value = '  -3.27%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '55.37'
$ws.Range('E37').Value = '  +0.37%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '499.51'
$ws.Range('E38').Value = '  -5.05%  '
$ws.Range('D39').Value = '0.0₃0770'
$ws.Range('E39').Value = '  +11.89%  '
$ws.Range('E40').Value = '  -2.81%  '
$ws.Range('E41').Value = '  +0.87%  '
$ws.Range('E42').Value = '  +3.61%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '8.71'
$ws.Range('E43').Value = '  -2.49%  '
$ws.Range('E44').Value = '  -3.29%  '
$ws.Range('D45').Value = '2.920.29'
$ws.Range('E45').Value = '  +0.48%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.44'
$ws.Range('E46').Value = '  -1.84%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '27.97'
$ws.Range('E47').Value = '  -2.75%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '2.39'
$ws.Range('E48').Value = '  +1.34%  '
$ws.Range('E50').Value = '  -0.84%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '121.30'
$ws.Range('E51').Value = '  -0.18%  '
